$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data per the scheduled data refresh.
# Force text format so numeric-looking strings (prices, percentages)
# keep their exact original formatting instead of being parsed as numbers.
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '28.376.61'
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = '  -0.40%  '
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.565.73'
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '  -0.05%  '
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = '  -0.06%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '211.02'
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.489'
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '  -0.68%  '
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '  -0.03%  '
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = '  -3.60%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '23.61'
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '  -1.86%  '
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '  -1.34%  '
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '  -0.63%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '1.789.01'
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '  -0.08%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '1.562.82'
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = '  -0.23%  '
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '  -0.23%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '28.358.11'
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '  -0.46%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '0.513'
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '  -1.28%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '60.60'
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '  -2.63%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '228.20'
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '  -0.11%  '
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '  +0.30%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '0.0₃0680'
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = '  -1.89%  '
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '  -0.05%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '3.94'
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '8.94'
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '  -1.97%  '
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '  -1.39%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '150.33'
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '  -0.26%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '14.89'
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = '  -0.72%  '
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = '  +0.23%  '
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = '  -2.00%  '
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = '  -0.03%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.0476'
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = '  +1.96%  '
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = '  -3.93%  '
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = '  -1.04%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '3.08'
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = '  +0.03%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '1.386.13'
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = '  -0.45%  '
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = '  +1.86%  '
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = '  -3.23%  '
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = '  -0.47%  '
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '  +2.41%  '
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '  -2.04%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.520'
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '  -3.10%  '
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '  +3.03%  '
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '  -0.05%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.785'
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '  -0.24%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.0469'
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = '  -2.23%  '
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '  -2.86%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '62.25'
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = '  -0.92%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.918'
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = '  -5.71%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '1.701.75'
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '  -0.24%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '85.51'
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = '  -0.49%  '
$c = $ws.Range("B51")
$c.NumberFormat = "@"
$c.Value = 'Cronos'
$c = $ws.Range("C51")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.0514'
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = '  -2.20%  '
